# "Generate Report for Handback"
# Updates the localization-status workbook to reflect that the de-de and
# zh-cn handback packages have come back from the translators:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (Overview + both per-locale sheets share this value).
#   - Each locale sheet's "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get populated now that a handback
#     has actually occurred.
#   - A hyperlink (to the source .md doc) is added on the new
#     "Latest Target File" cell, matching the existing one on column A.

$wb = $excel.ActiveWorkbook

$srcMdName  = "e9de683e-e0e0-4cc4-9f69-9c4bce5cc8e8.md"
$srcMdUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/128d41c1fac76b3253e92ed3493488cdfeb3552b/e2e/e9de683e-e0e0-4cc4-9f69-9c4bce5cc8e8.md"
$statusText = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status column updates (Overview + both locale sheets) ---
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$zhcn.Range("C2").Value = $statusText
$dede.Range("C2").Value = $statusText

# --- zh-cn: fill in handback columns ---
$zhcnXlf = $zhcn.Range("G2").Value2
$zhcn.Range("I2").Value = $srcMdName
$zhcn.Range("J2").Value = $zhcnXlf
$zhcn.Range("K2").Value = "2016-08-31 21:14:32"

$h = $zhcn.Hyperlinks.Add($zhcn.Range("I2"), $srcMdUrl, "", "", $srcMdName)

# --- de-de: fill in handback columns ---
$dedeXlf = $dede.Range("G2").Value2
$dede.Range("I2").Value = $srcMdName
$dede.Range("J2").Value = $dedeXlf
$dede.Range("K2").Value = "2016-08-31 21:14:40"

$h2 = $dede.Hyperlinks.Add($dede.Range("I2"), $srcMdUrl, "", "", $srcMdName)

# --- Column width adjustments (status / target-file / handback-file columns
#     widened to fit the longer content now held in them) ---
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(9).ColumnWidth = 39.2
$zhcn.Columns.Item(10).ColumnWidth = 39.2

$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(9).ColumnWidth = 39.2
$dede.Columns.Item(10).ColumnWidth = 39.2
